$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Helper: apply the bold / centered / wrapped "Total row" look used on the
# new Repayment-schedule total row.
# ---------------------------------------------------------------------------
function ApplyTotalFormat($rng) {
    $rng.Font.Bold = $true
    $rng.HorizontalAlignment = -4108   # xlCenter
    $rng.VerticalAlignment = -4108     # xlCenter
    $rng.WrapText = $true
}

# ---------------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("F2").Value = 0
$wsSummary.Range("A3").Value = 211.19
$wsSummary.Range("E3").Value = 114.2

# Refresh the sheet's remembered selection (A7:XFD13)
$wsSummary.Select()
$wsSummary.Range("A7:XFD13").Select()

# ---------------------------------------------------------------------------
# Repayment schedule sheet
# ---------------------------------------------------------------------------
$wsRepay = $wb.Worksheets.Item("Repayment schedule")

# -- Updated figures in the existing rows -----------------------------------
$wsRepay.Range("F4").Value = 921.65
$wsRepay.Range("G4").Value = 3211.57

$wsRepay.Range("F5").Value = 932.09
$wsRepay.Range("G5").Value = 2279.48

# H5 used to hold a formula; replace it with a plain number and drop the
# italic/number-format it inherited, falling back to the plain "data" style.
$wsRepay.Range("H5").ClearFormats()
$wsRepay.Range("H5").Value = 31.68
$wsRepay.Range("H5").VerticalAlignment = -4108
$wsRepay.Range("H5").WrapText = $true

$wsRepay.Range("F6").Value = 940.54
$wsRepay.Range("G6").Value = 1338.94
$wsRepay.Range("H6").Value = 23.23

$wsRepay.Range("F7").Value = 950.56
$wsRepay.Range("G7").Value = 388.38
$wsRepay.Range("H7").Value = 13.21

# D8/E8 were missing from the original sheet; give them the same empty,
# vertically centred + wrapped look as the rest of the data rows.
$wsRepay.Range("D8").VerticalAlignment = -4108
$wsRepay.Range("D8").WrapText = $true
$wsRepay.Range("E8").VerticalAlignment = -4108
$wsRepay.Range("E8").WrapText = $true

$wsRepay.Range("F8").Value = 388.38
$wsRepay.Range("H8").Value = 3.96
$wsRepay.Range("K8").Value = 392.34
$wsRepay.Range("P8").Value = 392.34

# -- New "Total" row (row 9) -------------------------------------------------
$wsRepay.Range("B9").Value = 212
$wsRepay.Range("C9").Value = "Total"
$wsRepay.Range("F9").Value = 5000
$wsRepay.Range("H9").Value = 211.19
$wsRepay.Range("I9").Value = 0
$wsRepay.Range("J9").Value = 0
$wsRepay.Range("K9").Value = 5211.1899999999996
$wsRepay.Range("L9").Value = 963.77
$wsRepay.Range("M9").Value = 0
$wsRepay.Range("N9").Value = 0
$wsRepay.Range("O9").Value = 0
$wsRepay.Range("P9").Value = 4247.42

# General (no special number format) cells of the total row
ApplyTotalFormat($wsRepay.Range("A9:B9"))
ApplyTotalFormat($wsRepay.Range("G9:J9"))
ApplyTotalFormat($wsRepay.Range("L9:O9"))

# Currency-ish (#,##0.00) cells
$rK9 = $wsRepay.Range("K9")
$rK9.NumberFormat = "#,##0.00"
ApplyTotalFormat($rK9)

$rP9 = $wsRepay.Range("P9")
$rP9.NumberFormat = "#,##0.00"
ApplyTotalFormat($rP9)

# Integer (#,##0) disbursed-amount cell
$rF9 = $wsRepay.Range("F9")
$rF9.NumberFormat = "#,##0"
ApplyTotalFormat($rF9)

# "Total" label, merged across C9:E9
$rC9 = $wsRepay.Range("C9:E9")
ApplyTotalFormat($rC9)
$rC9.Merge()

# Refresh the sheet's remembered selection (O2:O9)
$wsRepay.Select()
$wsRepay.Range("O2:O9").Select()

# ---------------------------------------------------------------------------
# Transactions sheet
# ---------------------------------------------------------------------------
$wsTrans = $wb.Worksheets.Item("Transactions")
$wsTrans.Range("A2").Value = 34
$wsTrans.Range("A3").Value = 32

# This is the sheet that ends up active/selected once the workbook is saved.
$wsTrans.Select()
$wsTrans.Range("A2:XFD4").Select()
